$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append row 8 with the new test-mail entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Ik wil mijn bestelling annuleren"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #5: Ik wil mijn bestelling annuleren"
$logs.Range("D8").Value = "Retour / Terugbetaling"
$logs.Range("E8").Value = "Beste klant,`nBedankt voor uw bericht. Om uw bestelling te annuleren, hebben we wat extra informatie nodig. Zou u zo vriendelijk willen zijn om uw bestelnummer met ons te delen? Op die manier kunnen we uw verzoek snel verwerken.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F8").Value = "2025-06-29 14:08:25"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"

# Extend the conditional-formatting ranges (D/G/H/I) from row 7 to row 8.
# Modifying any one rule's AppliesTo range updates the whole sqref group it
# belongs to, so one call per column group is enough.
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))

# --- Sheet "Dashboard": append new category row + extend chart series ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Retour / Terugbetaling"
$dash.Range("B5").Value = 1

$co = $dash.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$5,'Dashboard'!`$B`$2:`$B`$5,1)"
